function Set-TextValue($range, $value) {
    # Preserve the cell's existing style; force a Text number format so
    # numeric-looking strings (e.g. "300.42") are kept as literal text
    # instead of being parsed into a Double by Excel, then restore the
    # original style so no stray formatting is introduced.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "42.664.76"
Set-TextValue $ws.Range("E2") "  +0.67%  "
Set-TextValue $ws.Range("D3") "2.288.38"
Set-TextValue $ws.Range("E3") "  +0.36%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.08%  "
Set-TextValue $ws.Range("D5") "300.42"
Set-TextValue $ws.Range("E5") "  +0.04%  "
Set-TextValue $ws.Range("D6") "99.04"
Set-TextValue $ws.Range("E6") "  +2.73%  "
Set-TextValue $ws.Range("D7") "0.500"
Set-TextValue $ws.Range("E7") "  +0.63%  "
Set-TextValue $ws.Range("E8") "  -0.10%  "
Set-TextValue $ws.Range("E9") "  +3.95%  "
Set-TextValue $ws.Range("D10") "35.81"
Set-TextValue $ws.Range("E10") "  +7.42%  "
Set-TextValue $ws.Range("D11") "0.0789"
Set-TextValue $ws.Range("E11") "  -0.22%  "
Set-TextValue $ws.Range("E12") "  +2.15%  "
Set-TextValue $ws.Range("D13") "17.79"
Set-TextValue $ws.Range("E13") "  +11.16%  "
Set-TextValue $ws.Range("D14") "6.81"
Set-TextValue $ws.Range("E14") "  +1.44%  "
Set-TextValue $ws.Range("D15") "2.645.31"
Set-TextValue $ws.Range("E15") "  +0.45%  "
Set-TextValue $ws.Range("D16") "2.283.09"
Set-TextValue $ws.Range("E16") "  +0.14%  "
Set-TextValue $ws.Range("E17") "  +0.78%  "
Set-TextValue $ws.Range("D18") "42.578.51"
Set-TextValue $ws.Range("E18") "  +0.58%  "
Set-TextValue $ws.Range("D19") "12.40"
Set-TextValue $ws.Range("E19") "  +5.88%  "
Set-TextValue $ws.Range("E20") "  +3.05%  "
Set-TextValue $ws.Range("E21") "  +0.36%  "
Set-TextValue $ws.Range("D22") "67.77"
Set-TextValue $ws.Range("E22") "  +1.89%  "
Set-TextValue $ws.Range("D23") "235.34"
Set-TextValue $ws.Range("E23") "  -0.09%  "
Set-TextValue $ws.Range("D24") "2.21"
Set-TextValue $ws.Range("E24") "  +12.75%  "
Set-TextValue $ws.Range("E25") "  -0.07%  "
Set-TextValue $ws.Range("E26") "  -0.47%  "
Set-TextValue $ws.Range("D27") "24.54"
Set-TextValue $ws.Range("E27") "  +3.07%  "
Set-TextValue $ws.Range("B28") "Monero"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "168.23"
Set-TextValue $ws.Range("E28") "  +0.39%  "
Set-TextValue $ws.Range("B29") "InjectiveProtocol"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D29") "34.30"
Set-TextValue $ws.Range("E29") "  +1.67%  "
Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.06"
Set-TextValue $ws.Range("E30") "  -0.40%  "
Set-TextValue $ws.Range("D31") "9.15"
Set-TextValue $ws.Range("E31") "  +0.00%  "
Set-TextValue $ws.Range("E32") "  -0.13%  "
Set-TextValue $ws.Range("E33") "  +0.54%  "
Set-TextValue $ws.Range("D34") "17.54"
Set-TextValue $ws.Range("E34") "  +4.19%  "
Set-TextValue $ws.Range("E35") "  -1.93%  "
Set-TextValue $ws.Range("E36") "  +2.85%  "
Set-TextValue $ws.Range("E37") "  -0.86%  "
Set-TextValue $ws.Range("E38") "  +1.26%  "
Set-TextValue $ws.Range("E39") "  +2.08%  "
Set-TextValue $ws.Range("D40") "1.77"
Set-TextValue $ws.Range("E40") "  +2.10%  "
Set-TextValue $ws.Range("E41") "  +0.22%  "
Set-TextValue $ws.Range("D42") "1.989.19"
Set-TextValue $ws.Range("E42") "  +1.63%  "
Set-TextValue $ws.Range("E43") "  +2.65%  "
Set-TextValue $ws.Range("E44") "  -1.92%  "
Set-TextValue $ws.Range("D45") "10.10"
Set-TextValue $ws.Range("E45") "  +4.78%  "
Set-TextValue $ws.Range("E46") "  +2.53%  "
Set-TextValue $ws.Range("D47") "17.45"
Set-TextValue $ws.Range("E47") "  -0.06%  "
Set-TextValue $ws.Range("D48") "55.42"
Set-TextValue $ws.Range("E48") "  +5.75%  "
Set-TextValue $ws.Range("D49") "2.513.10"
Set-TextValue $ws.Range("E49") "  +0.37%  "
Set-TextValue $ws.Range("E50") "  +2.76%  "
Set-TextValue $ws.Range("E51") "  -0.69%  "
